$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Femacal de La Calera - Cereza, "Early Burlat")
# is inserted as row 163; every following record shifts down by one row
# (old row 196 becomes the new row 197).
$ws.Rows(163).Insert()

$ws.Range("A163").Value = 3
$ws.Range("B163").Value = "Femacal de La Calera"
$ws.Range("C163").Value = "Coquimbo"
$ws.Range("D163").Value = 44504
$ws.Range("E163").Value = 5
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100103
$ws.Range("H163").Value = "Frutos de hueso (carozo)"
$ws.Range("I163").Value = 100103001
$ws.Range("J163").Value = "Cereza"
$ws.Range("K163").Value = "Early Burlat"
$ws.Range("L163").Value = "Primera"
$ws.Range("M163").Value = 65
$ws.Range("N163").Value = 50000
$ws.Range("O163").Value = 50000
$ws.Range("P163").Value = 50000
$ws.Range("Q163").Value = "$/bandeja 10 kilos"
$ws.Range("R163").Value = "Provincia de Curicó"
$ws.Range("S163").Value = 5000
$ws.Range("T163").Value = 10
